$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("E3").Value = Get-Date -Year 2025 -Month 6 -Day 6 -Hour 11 -Minute 0 -Second 0
$ws.Range("E3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E4").Value = Get-Date -Year 2025 -Month 6 -Day 12 -Hour 9 -Minute 0 -Second 0
$ws.Range("E5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E5").Value = Get-Date -Year 2025 -Month 6 -Day 12 -Hour 9 -Minute 0 -Second 0
$ws.Range("E6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E6").Value = Get-Date -Year 2025 -Month 6 -Day 12 -Hour 9 -Minute 0 -Second 0
$ws.Range("E7").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E7").Value = Get-Date -Year 2025 -Month 6 -Day 12 -Hour 9 -Minute 0 -Second 0
$ws.Range("E8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E8").Value = Get-Date -Year 2025 -Month 6 -Day 10 -Hour 11 -Minute 0 -Second 0
$ws.Range("E12").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E12").Value = Get-Date -Year 2025 -Month 6 -Day 12 -Hour 9 -Minute 0 -Second 0
$ws.Range("E13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E13").Value = Get-Date -Year 2025 -Month 6 -Day 12 -Hour 9 -Minute 0 -Second 0
$ws.Range("E14").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E14").Value = Get-Date -Year 2025 -Month 6 -Day 12 -Hour 9 -Minute 0 -Second 0
$ws.Range("E15").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E15").Value = Get-Date -Year 2025 -Month 6 -Day 12 -Hour 9 -Minute 0 -Second 0
$ws.Range("E16").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E16").Value = Get-Date -Year 2025 -Month 6 -Day 12 -Hour 9 -Minute 0 -Second 0
$ws.Range("E17").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E17").Value = Get-Date -Year 2025 -Month 6 -Day 12 -Hour 10 -Minute 0 -Second 0
$ws.Range("E18").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E18").Value = Get-Date -Year 2025 -Month 6 -Day 10 -Hour 10 -Minute 0 -Second 0
$ws.Range("E20").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E20").Value = Get-Date -Year 2025 -Month 6 -Day 10 -Hour 10 -Minute 0 -Second 0
$ws.Range("E34").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E34").Value = Get-Date -Year 2025 -Month 6 -Day 6 -Hour 10 -Minute 0 -Second 0
$ws.Range("E35").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E35").Value = Get-Date -Year 2025 -Month 6 -Day 10 -Hour 11 -Minute 0 -Second 0
$ws.Range("E37").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E37").Value = Get-Date -Year 2025 -Month 6 -Day 18 -Hour 10 -Minute 0 -Second 0
$ws.Range("E39").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E39").Value = Get-Date -Year 2025 -Month 6 -Day 4 -Hour 11 -Minute 0 -Second 0
$ws.Range("E40").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E40").Value = Get-Date -Year 2025 -Month 6 -Day 4 -Hour 10 -Minute 0 -Second 0
$ws.Range("E41").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E41").Value = Get-Date -Year 2025 -Month 6 -Day 6 -Hour 11 -Minute 0 -Second 0
$ws.Range("E42").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E42").Value = Get-Date -Year 2025 -Month 6 -Day 13 -Hour 12 -Minute 0 -Second 0
$ws.Range("E43").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E43").Value = Get-Date -Year 2025 -Month 6 -Day 18 -Hour 10 -Minute 0 -Second 0
$ws.Range("E48").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E48").Value = Get-Date -Year 2025 -Month 6 -Day 2 -Hour 10 -Minute 0 -Second 0
$ws.Range("E49").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E49").Value = Get-Date -Year 2025 -Month 6 -Day 11 -Hour 10 -Minute 0 -Second 0
$ws.Range("E50").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E50").Value = Get-Date -Year 2025 -Month 6 -Day 2 -Hour 11 -Minute 0 -Second 0
$ws.Range("E51").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E51").Value = Get-Date -Year 2025 -Month 6 -Day 6 -Hour 10 -Minute 0 -Second 0
$ws.Range("E52").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E52").Value = Get-Date -Year 2025 -Month 6 -Day 5 -Hour 10 -Minute 0 -Second 0
$ws.Range("E53").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E53").Value = Get-Date -Year 2025 -Month 6 -Day 6 -Hour 10 -Minute 0 -Second 0
$ws.Range("E54").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E54").Value = Get-Date -Year 2025 -Month 6 -Day 10 -Hour 10 -Minute 0 -Second 0
$ws.Range("E55").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E55").Value = Get-Date -Year 2025 -Month 6 -Day 11 -Hour 11 -Minute 0 -Second 0
$ws.Range("E56").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E56").Value = Get-Date -Year 2025 -Month 6 -Day 10 -Hour 11 -Minute 0 -Second 0
$ws.Range("E57").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E57").Value = Get-Date -Year 2025 -Month 6 -Day 5 -Hour 11 -Minute 0 -Second 0
$ws.Range("E58").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E58").Value = Get-Date -Year 2025 -Month 6 -Day 5 -Hour 10 -Minute 0 -Second 0
$ws.Range("E61").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E61").Value = Get-Date -Year 2025 -Month 6 -Day 3 -Hour 10 -Minute 0 -Second 0
$ws.Range("E66").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E66").Value = Get-Date -Year 2025 -Month 6 -Day 3 -Hour 11 -Minute 0 -Second 0
$ws.Range("E67").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E67").Value = Get-Date -Year 2025 -Month 6 -Day 10 -Hour 11 -Minute 0 -Second 0
$ws.Range("E69").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E69").Value = Get-Date -Year 2025 -Month 5 -Day 30 -Hour 10 -Minute 0 -Second 0
$ws.Range("E70").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E70").Value = Get-Date -Year 2025 -Month 5 -Day 29 -Hour 10 -Minute 0 -Second 0
$ws.Range("E72").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E72").Value = Get-Date -Year 2025 -Month 5 -Day 30 -Hour 11 -Minute 0 -Second 0
$ws.Range("E76").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E76").Value = Get-Date -Year 2025 -Month 6 -Day 2 -Hour 11 -Minute 0 -Second 0
$ws.Range("E78").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E78").Value = Get-Date -Year 2025 -Month 5 -Day 30 -Hour 10 -Minute 0 -Second 0
$ws.Range("E81").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E81").Value = Get-Date -Year 2025 -Month 6 -Day 3 -Hour 11 -Minute 0 -Second 0
$ws.Range("E82").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E82").Value = Get-Date -Year 2025 -Month 5 -Day 29 -Hour 9 -Minute 30 -Second 0
$ws.Range("E83").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E83").Value = Get-Date -Year 2025 -Month 5 -Day 29 -Hour 11 -Minute 0 -Second 0
$ws.Range("E84").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E84").Value = Get-Date -Year 2025 -Month 5 -Day 29 -Hour 11 -Minute 0 -Second 0
$ws.Range("E85").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E85").Value = Get-Date -Year 2025 -Month 5 -Day 30 -Hour 11 -Minute 0 -Second 0
$ws.Range("E86").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E86").Value = Get-Date -Year 2025 -Month 5 -Day 30 -Hour 10 -Minute 0 -Second 0
$ws.Range("E90").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E90").Value = Get-Date -Year 2025 -Month 5 -Day 26 -Hour 10 -Minute 0 -Second 0
$ws.Range("E91").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E91").Value = Get-Date -Year 2025 -Month 5 -Day 19 -Hour 11 -Minute 0 -Second 0
$ws.Range("E92").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E92").Value = Get-Date -Year 2025 -Month 5 -Day 26 -Hour 9 -Minute 30 -Second 0
$ws.Range("E93").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E93").Value = Get-Date -Year 2025 -Month 5 -Day 23 -Hour 11 -Minute 0 -Second 0
$ws.Range("E94").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E94").Value = Get-Date -Year 2025 -Month 5 -Day 27 -Hour 10 -Minute 0 -Second 0
$ws.Range("E95").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E95").Value = Get-Date -Year 2025 -Month 5 -Day 26 -Hour 10 -Minute 0 -Second 0
$ws.Range("E96").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E96").Value = Get-Date -Year 2024 -Month 5 -Day 23 -Hour 11 -Minute 0 -Second 0
$ws.Range("E97").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E97").Value = Get-Date -Year 2025 -Month 5 -Day 19 -Hour 13 -Minute 0 -Second 0
$ws.Range("E98").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E98").Value = Get-Date -Year 2025 -Month 5 -Day 22 -Hour 10 -Minute 0 -Second 0
$ws.Range("E99").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E99").Value = Get-Date -Year 2025 -Month 5 -Day 21 -Hour 0 -Minute 0 -Second 0
$ws.Range("E100").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E100").Value = Get-Date -Year 2025 -Month 5 -Day 22 -Hour 10 -Minute 12 -Second 0
$ws.Range("E101").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E101").Value = Get-Date -Year 2025 -Month 5 -Day 19 -Hour 11 -Minute 0 -Second 0
$ws.Range("E103").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E103").Value = Get-Date -Year 2025 -Month 5 -Day 21 -Hour 11 -Minute 0 -Second 0
$ws.Range("E104").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E104").Value = Get-Date -Year 2025 -Month 5 -Day 16 -Hour 11 -Minute 0 -Second 0
$ws.Range("E105").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E105").Value = Get-Date -Year 2025 -Month 5 -Day 13 -Hour 10 -Minute 0 -Second 0
$ws.Range("E106").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E106").Value = Get-Date -Year 2025 -Month 5 -Day 19 -Hour 11 -Minute 0 -Second 0
$ws.Range("E107").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E107").Value = Get-Date -Year 2025 -Month 5 -Day 14 -Hour 11 -Minute 0 -Second 0
$ws.Range("E108").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E108").Value = Get-Date -Year 2025 -Month 5 -Day 6 -Hour 11 -Minute 0 -Second 0
$ws.Range("E109").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E109").Value = Get-Date -Year 2025 -Month 4 -Day 25 -Hour 11 -Minute 0 -Second 0
$ws.Range("E110").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E110").Value = Get-Date -Year 2025 -Month 5 -Day 7 -Hour 10 -Minute 0 -Second 0
$ws.Range("E111").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E111").Value = Get-Date -Year 2025 -Month 5 -Day 6 -Hour 11 -Minute 0 -Second 0
$ws.Range("E112").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E112").Value = Get-Date -Year 2025 -Month 4 -Day 17 -Hour 10 -Minute 0 -Second 0
$ws.Range("E113").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E113").Value = Get-Date -Year 2025 -Month 4 -Day 10 -Hour 10 -Minute 0 -Second 0
$ws.Range("E114").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E114").Value = Get-Date -Year 2025 -Month 3 -Day 27 -Hour 10 -Minute 0 -Second 0
